# "update of results and scripts. Anonimyzed fedcore"
#
# For every sheet:
#   1. Rename the "fedcore" column header(s) in row 2 to "approach".
#   2. Give the blank spacer cells that sit inside each merged group header
#      (e.g. C1/D1 under the merged B1:D1 band) a light top+bottom border
#      (and, for the right-most spacer of the band, also a right border)
#      instead of the heavy all-round thin border they inherited from the
#      default header style.
# Then, on the computational_comparison sheet only, drop the stray empty
# cell G5.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # --- 1. "fedcore" -> "approach" ---------------------------------------
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $used.Cells.Item($r, $c)
            if ($cell.Value2 -eq "fedcore") {
                $cell.Value = "approach"
            }
        }
    }

    # --- 2. Re-border the spacer cells of each merged header band ---------
    # Every merged area in row 1 (e.g. B1:D1, E1:G1, ...) is a 3-column
    # group: the first column carries the "0"/"1" index value, and the two
    # columns to its right are blank spacers that get the lighter border.
    $used = $ws.UsedRange
    $seenGroups = @{}
    for ($c = 1; $c -le $used.Columns.Count; $c++) {
        $headCell = $ws.Cells.Item(1, $c)
        if (-not $headCell.MergeCells) { continue }

        $area = $headCell.MergeArea
        $addr = $area.Address()
        if ($seenGroups.ContainsKey($addr)) { continue }
        $seenGroups[$addr] = $true
        if ($area.Columns.Count -lt 3) { continue }

        $leftCol = $area.Column + 1
        $rightCol = $area.Column + $area.Columns.Count - 1

        $leftCell = $ws.Cells.Item(1, $leftCol)
        $rightCell = $ws.Cells.Item(1, $rightCol)

        # Apply the shared top+bottom border to both spacer cells together
        # (as one contiguous range) so they end up sharing a single style.
        $pairRange = $ws.Range($leftCell, $rightCell)
        $pairRange.ClearFormats()
        $pairRange.Borders.Item(8).LineStyle = 1   # xlEdgeTop
        $pairRange.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

        # The right-most spacer additionally gets a right border.
        $rightCell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    }
}

# --- 3. Drop the stray empty cell G5 on computational_comparison ----------
$ws2 = $wb.Worksheets.Item("computational_comparison")
$ws2.Range("G5").ClearContents()
